# Actualización desde MV -datos-
# Adds a new quarterly row (01-04-2021) to the debt-by-instrument table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 79

# Column A holds a date-like label that must be stored as literal text
# (matching the existing "01-01-2002".."01-01-2021" entries), not an
# auto-converted Excel date serial. Temporarily force a Text number
# format so the assignment isn't reinterpreted as a date, then clear
# the formatting back to the sheet's default so the cell matches its
# siblings (no explicit style override).
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "01-04-2021"
$cellA.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 97831
$ws.Cells.Item($newRow, 3).Value = 96245
$ws.Cells.Item($newRow, 4).Value = 19
$ws.Cells.Item($newRow, 5).Value = 1567
$ws.Cells.Item($newRow, 6).Value = 69174
$ws.Cells.Item($newRow, 7).Value = 69155
$ws.Cells.Item($newRow, 8).Value = 19
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 28657
$ws.Cells.Item($newRow, 11).Value = 27090
$ws.Cells.Item($newRow, 12).Value = 1567
